# Update Price (D) and Volume(1h) (E) text columns to reflect the latest
# cryptocurrency snapshot, preserving each value's original text formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.411.50"
$ws.Cells.Item(3, 4).Value = "1.671.58"
$ws.Cells.Item(3, 5).Value = "  +1.12%  "
$ws.Cells.Item(4, 5).Value = "  +0.52%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "221.24"
$ws.Cells.Item(5, 5).Value = "  +1.66%  "
$ws.Cells.Item(6, 5).Value = "  +0.77%  "
$ws.Cells.Item(7, 5).Value = "  +0.52%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.2663"
$ws.Cells.Item(8, 5).Value = "  +1.61%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.06380"
$ws.Cells.Item(9, 5).Value = "  +0.94%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "20.87"
$ws.Cells.Item(10, 5).Value = "  +2.34%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.07864"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "4.530"
$ws.Cells.Item(12, 5).Value = "  +0.30%  "
$ws.Cells.Item(13, 4).Value = "1.675.89"
$ws.Cells.Item(13, 5).Value = "  +0.26%  "
$ws.Cells.Item(14, 4).Value = "1.901.34"
$ws.Cells.Item(14, 5).Value = "  +1.08%  "
$ws.Cells.Item(15, 5).Value = "  +2.30%  "
$ws.Cells.Item(16, 4).Value = "0.0₅8193"
$ws.Cells.Item(16, 5).Value = "  +0.27%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "66.17"
$ws.Cells.Item(17, 5).Value = "  +1.28%  "
$ws.Cells.Item(18, 4).Value = "26.424.25"
$ws.Cells.Item(18, 5).Value = "  +1.05%  "
$ws.Cells.Item(19, 5).Value = "  +0.53%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "4.722"
$ws.Cells.Item(20, 5).Value = "  +2.75%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "197.98"
$ws.Cells.Item(21, 5).Value = "  +3.65%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "10.31"
$ws.Cells.Item(22, 5).Value = "  +2.26%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "6.073"
$ws.Cells.Item(23, 5).Value = "  +1.22%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "1.011"
$ws.Cells.Item(24, 5).Value = "  +0.42%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "145.88"
$ws.Cells.Item(25, 5).Value = "  +0.40%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.1226"
$ws.Cells.Item(26, 5).Value = "  +0.27%  "
$ws.Cells.Item(27, 5).Value = "  +0.77%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "16.25"
$ws.Cells.Item(28, 5).Value = "  +1.69%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.509"
$ws.Cells.Item(29, 5).Value = "  +2.70%  "
$ws.Cells.Item(30, 5).Value = "  +3.38%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.288"
$ws.Cells.Item(31, 5).Value = "  +1.24%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.564"
$ws.Cells.Item(32, 5).Value = "  +0.54%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "3.334"
$ws.Cells.Item(33, 5).Value = "  +2.29%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.610"
$ws.Cells.Item(34, 5).Value = "  +1.44%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.9695"
$ws.Cells.Item(35, 5).Value = "  +2.35%  "
$ws.Cells.Item(36, 5).Value = "  +1.22%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.5839"
$ws.Cells.Item(38, 5).Value = "  +2.00%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.01617"
$ws.Cells.Item(39, 5).Value = "  +0.66%  "
$ws.Cells.Item(40, 4).Value = "1.080.56"
$ws.Cells.Item(40, 5).Value = "  +4.12%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "5.936"
$ws.Cells.Item(41, 5).Value = "  +2.43%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.8650"
$ws.Cells.Item(42, 5).Value = "  +1.54%  "
$ws.Cells.Item(43, 5).Value = "  +0.56%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "103.11"
$ws.Cells.Item(44, 5).Value = "  -0.73%  "
$ws.Cells.Item(45, 4).Value = "1.812.38"
$ws.Cells.Item(45, 5).Value = "  +1.05%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "58.52"
$ws.Cells.Item(46, 5).Value = "  +3.24%  "
$ws.Cells.Item(47, 4).Value = "0.0₈107"
$ws.Cells.Item(47, 5).Value = "  +3.91%  "
$ws.Cells.Item(48, 5).Value = "  +0.68%  "
$ws.Cells.Item(49, 5).Value = "  +1.44%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "8.034"
$ws.Cells.Item(50, 5).Value = "  +2.45%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.05161"
$ws.Cells.Item(51, 5).Value = "  +0.13%  "
